$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "Session" to "Neurology"
$ws.Name = "Neurology"

# Append the two new QR-scanner log rows (37 and 38)
$ws.Range("A37").Value = "'190333"
$ws.Range("B37").Value = "Neurology"
$ws.Range("C37").Value = "16/12/2025"
$ws.Range("D37").Value = "10:13:46"
$ws.Range("E37").Value = "Manual"
$ws.Range("F37").Value = "emp17.farah.a.youssef@gmail.com"

$ws.Range("A38").Value = "'191007"
$ws.Range("B38").Value = "Neurology"
$ws.Range("C38").Value = "16/12/2025"
$ws.Range("D38").Value = "10:16:24"
$ws.Range("E38").Value = "Scan"
$ws.Range("F38").Value = "emp17.farah.a.youssef@gmail.com"
